# "fixes + videos in BO + discord in BO"
#
# Sheet1 (Blad1)  = Water 2 TC FH build (Anatolia) - By Kosis
# Sheet2 (Feuil1) = Water FH build (Anatolia) - By Kosis
# Sheet3 (Feuil2) = Eco upgrades - By Looki

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- text fixes on the "Eco upgrades" sheet -------------------------------
$ws3.Range("B4").Value  = "Initial vils : 1 build oxcart (send it to gold) -> hunt / 2 wood  "
$ws3.Range("A15").Value = "remark :"

# --- add a youtube link under the "Water FH build" sheet ------------------
$ws2.Range("A28").Value = "https://www.youtube.com/watch?v=kdWf_uQ3xgw"
[void]$ws2.Range("A28").Select()

# --- fix the lingering selection on the "Eco upgrades" sheet --------------
[void]$ws3.Range("A15").Select()

# --- add a youtube/discord link under the "Water 2 TC FH build" sheet -----
# (done last so this sheet/selection ends up the active one on reopen)
$ws1.Range("A31").Value = "https://youtu.be/kdWf_uQ3xgw?si=DCtXz05A6b4IafTj&t=437"
[void]$ws1.Range("A31").Select()
